# Insert 4 new data rows for Lechuga / Macroferia Regional de Talca
# right after the existing "Marina" row for date 44299 (current row 1221),
# pushing the existing rows 1222..1309 down to 1226..1313.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("1222:1225").Insert()

# New row 1222: Conconina(o)
$r = 1222
$ws.Cells.Item($r,1).Value  = 5
$ws.Cells.Item($r,2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item($r,3).Value  = "Maule"
$ws.Cells.Item($r,4).Value  = 44714
$ws.Cells.Item($r,5).Value  = 7
$ws.Cells.Item($r,6).Value  = 100112033
$ws.Cells.Item($r,7).Value  = "Lechuga"
$ws.Cells.Item($r,8).Value  = "Conconina(o)"
$ws.Cells.Item($r,9).Value  = "Primera"
$ws.Cells.Item($r,10).Value = 500
$ws.Cells.Item($r,11).Value = 4000
$ws.Cells.Item($r,12).Value = 4000
$ws.Cells.Item($r,13).Value = 4000
$ws.Cells.Item($r,14).Value = "`$/caja 10 unidades"
$ws.Cells.Item($r,15).Value = "Región del Maule"
$ws.Cells.Item($r,16).Value = 400
$ws.Cells.Item($r,17).Value = 10
$ws.Cells.Item($r,18).Value = "Hortaliza"

# New row 1223: Escarola
$r = 1223
$ws.Cells.Item($r,1).Value  = 5
$ws.Cells.Item($r,2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item($r,3).Value  = "Maule"
$ws.Cells.Item($r,4).Value  = 44714
$ws.Cells.Item($r,5).Value  = 7
$ws.Cells.Item($r,6).Value  = 100112033
$ws.Cells.Item($r,7).Value  = "Lechuga"
$ws.Cells.Item($r,8).Value  = "Escarola"
$ws.Cells.Item($r,9).Value  = "Primera"
$ws.Cells.Item($r,10).Value = 700
$ws.Cells.Item($r,11).Value = 6500
$ws.Cells.Item($r,12).Value = 6500
$ws.Cells.Item($r,13).Value = 6500
$ws.Cells.Item($r,14).Value = "`$/caja 15 unidades"
$ws.Cells.Item($r,15).Value = "Provincia del Elquí"
$ws.Cells.Item($r,16).Value = 433
$ws.Cells.Item($r,17).Value = 15
$ws.Cells.Item($r,18).Value = "Hortaliza"

# New row 1224: Española
$r = 1224
$ws.Cells.Item($r,1).Value  = 5
$ws.Cells.Item($r,2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item($r,3).Value  = "Maule"
$ws.Cells.Item($r,4).Value  = 44714
$ws.Cells.Item($r,5).Value  = 7
$ws.Cells.Item($r,6).Value  = 100112033
$ws.Cells.Item($r,7).Value  = "Lechuga"
$ws.Cells.Item($r,8).Value  = "Española"
$ws.Cells.Item($r,9).Value  = "Primera"
$ws.Cells.Item($r,10).Value = 400
$ws.Cells.Item($r,11).Value = 4000
$ws.Cells.Item($r,12).Value = 4000
$ws.Cells.Item($r,13).Value = 4000
$ws.Cells.Item($r,14).Value = "`$/caja 18 unidades"
$ws.Cells.Item($r,15).Value = "Región del Maule"
$ws.Cells.Item($r,16).Value = 222
$ws.Cells.Item($r,17).Value = 18
$ws.Cells.Item($r,18).Value = "Hortaliza"

# New row 1225: Marina
$r = 1225
$ws.Cells.Item($r,1).Value  = 5
$ws.Cells.Item($r,2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item($r,3).Value  = "Maule"
$ws.Cells.Item($r,4).Value  = 44714
$ws.Cells.Item($r,5).Value  = 7
$ws.Cells.Item($r,6).Value  = 100112033
$ws.Cells.Item($r,7).Value  = "Lechuga"
$ws.Cells.Item($r,8).Value  = "Marina"
$ws.Cells.Item($r,9).Value  = "Primera"
$ws.Cells.Item($r,10).Value = 500
$ws.Cells.Item($r,11).Value = 4000
$ws.Cells.Item($r,12).Value = 4000
$ws.Cells.Item($r,13).Value = 4000
$ws.Cells.Item($r,14).Value = "`$/caja 18 unidades"
$ws.Cells.Item($r,15).Value = "Región del Maule"
$ws.Cells.Item($r,16).Value = 222
$ws.Cells.Item($r,17).Value = 18
$ws.Cells.Item($r,18).Value = "Hortaliza"
